$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.025.26'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').Value = '3.091.61'
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.14'
$ws.Range('E5').Value = '  +8.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '619.38'
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('E7').Value = '  -12.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.358'
$ws.Range('E8').Value = '  -3.54%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').Value = '3.089.44'
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.715'
$ws.Range('E11').Value = '  -6.98%  '
$ws.Range('E12').Value = '  -3.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.03'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').Value = '89.833.77'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.37'
$ws.Range('E16').Value = '  -6.59%  '
$ws.Range('D17').Value = '3.653.76'
$ws.Range('E17').Value = '  -2.78%  '
$ws.Range('D18').Value = '3.069.44'
$ws.Range('E18').Value = '  -3.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.79'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000210'
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.81'
$ws.Range('E21').Value = '  -6.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '431.17'
$ws.Range('E22').Value = '  -9.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.44'
$ws.Range('E23').Value = '  +2.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.76'
$ws.Range('E24').Value = '  -4.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.57'
$ws.Range('E25').Value = '  -6.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '86.00'
$ws.Range('E26').Value = '  -11.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.75'
$ws.Range('E27').Value = '  -5.36%  '
$ws.Range('D28').Value = '3.280.55'
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.10'
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('E31').Value = '  +3.00%  '
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.57'
$ws.Range('E34').Value = '  -9.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.150'
$ws.Range('E35').Value = '  +4.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.71'
$ws.Range('E36').Value = '  +2.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '495.99'
$ws.Range('E37').Value = '  -5.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.04'
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.88'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.26'
$ws.Range('E40').Value = '  -4.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').Value = '  +54.45%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0864'
$ws.Range('E42').Value = '  -4.43%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.09'
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.398'
$ws.Range('E45').Value = '  -5.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.87'
$ws.Range('E46').Value = '  -6.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.676'
$ws.Range('E47').Value = '  -4.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '150.18'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('E49').Value = '  -2.16%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.31'
$ws.Range('E50').Value = '  -5.21%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.04%  '
